# Weekly update: prepend a new week's worth of "Apio" records (Primera y
# Segunda calidad) for Vega Central Mapocho de Santiago, pushing the
# existing records down by two rows. The two rows that fall off the
# bottom of the contiguous data block are re-appended as new rows at the
# end of the sheet (354 and 355), growing the used range from
# A1:R353 to A1:R355.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 249..353 down to 251..355, leaving two blank rows (249:250)
# ready for the new week's data.
$ws.Range("249:250").EntireRow.Insert()

# New row 249: Apio, Americana (o), Primera
$ws.Range("A249").Value = 9
$ws.Range("B249").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C249").Value = "Metropolitana"
$ws.Range("D249").Value = 44917
$ws.Range("E249").Value = 13
$ws.Range("F249").Value = 100112017
$ws.Range("G249").Value = "Apio"
$ws.Range("H249").Value = "Americana (o)"
$ws.Range("I249").Value = "Primera"
$ws.Range("J249").Value = 90
$ws.Range("K249").Value = 8000
$ws.Range("L249").Value = 9000
$ws.Range("M249").Value = 8500
$ws.Range("N249").Value = "$/docena de matas"
$ws.Range("O249").Value = "Región de Coquimbo"
$ws.Range("P249").Value = 1417
$ws.Range("Q249").Value = 6
$ws.Range("R249").Value = "Hortaliza"

# New row 250: Apio, Americana (o), Segunda
$ws.Range("A250").Value = 9
$ws.Range("B250").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C250").Value = "Metropolitana"
$ws.Range("D250").Value = 44917
$ws.Range("E250").Value = 13
$ws.Range("F250").Value = 100112017
$ws.Range("G250").Value = "Apio"
$ws.Range("H250").Value = "Americana (o)"
$ws.Range("I250").Value = "Segunda"
$ws.Range("J250").Value = 44
$ws.Range("K250").Value = 7000
$ws.Range("L250").Value = 7000
$ws.Range("M250").Value = 7000
$ws.Range("N250").Value = "$/docena de matas"
$ws.Range("O250").Value = "Región de Coquimbo"
$ws.Range("P250").Value = 1167
$ws.Range("Q250").Value = 6
$ws.Range("R250").Value = "Hortaliza"
